$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2 through 252 from 45192 to 45202
for ($r = 2; $r -le 252; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# Add new row 253 with data
$ws.Cells.Item(253, 1).Value = "A 45421-2023"
$ws.Cells.Item(253, 2).Value = 45194
$ws.Cells.Item(253, 3).Value = 45202
$ws.Cells.Item(253, 4).Value = "DALARNAS LÄN"
$ws.Cells.Item(253, 5).Value = "GAGNEF"
$ws.Cells.Item(253, 7).Value = 1.6
$ws.Cells.Item(253, 8).Value = 0
$ws.Cells.Item(253, 9).Value = 0
$ws.Cells.Item(253, 10).Value = 0
$ws.Cells.Item(253, 11).Value = 0
$ws.Cells.Item(253, 12).Value = 0
$ws.Cells.Item(253, 13).Value = 0
$ws.Cells.Item(253, 14).Value = 0
$ws.Cells.Item(253, 15).Value = 0
$ws.Cells.Item(253, 16).Value = 0
$ws.Cells.Item(253, 17).Value = 0

# Copy style (date format) from B252/C252 to B253/C253
$ws.Range("B252:C252").Copy()
$ws.Range("B253:C253").PasteSpecial(-4122)  # xlPasteFormats

# Copy style from R252 (wrap text) to R253
$ws.Range("R252").Copy()
$ws.Range("R253").PasteSpecial(-4122)

$ws.Range("R253").Value = ""

# Row 252 needs explicit custom row height like the rest of the sheet
$ws.Rows.Item(252).RowHeight = 15

$excel.CutCopyMode = $false
